$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 103 (shifts old rows 103-118 down to 104-119),
# preserving all existing formatting/values below.
$ws.Rows(103).Insert()

# Populate the newly inserted row 103 with the new weekly record.
# Values mirror the old row 103 (now row 104) except for the date (D) and
# volume (J), which carry the new week's figures.
$ws.Range("A103").Value = 5
$ws.Range("B103").Value = "Macroferia Regional de Talca"
$ws.Range("C103").Value = "Maule"
$ws.Range("D103").Value = 44474
$ws.Range("E103").Value = 7
$ws.Range("F103").Value = 100112017
$ws.Range("G103").Value = "Apio"
$ws.Range("H103").Value = "Americana (o)"
$ws.Range("I103").Value = "Primera"
$ws.Range("J103").Value = 800
$ws.Range("K103").Value = 7000
$ws.Range("L103").Value = 7000
$ws.Range("M103").Value = 7000
$ws.Range("N103").Value = "$/docena de matas"
$ws.Range("O103").Value = "Provincia del Elquí"
$ws.Range("P103").Value = 1167
$ws.Range("Q103").Value = 6
$ws.Range("R103").Value = "Hortaliza"
